$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.599.59"
$ws.Range("E2").Value = "  -0.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.726.27"
$ws.Range("E3").Value = "  -0.76%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "564.45"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.79"
$ws.Range("E6").Value = "  +1.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  -0.67%  "

$ws.Range("E9").Value = "  +0.35%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +4.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.63"
$ws.Range("E11").Value = "  +0.81%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.379"
$ws.Range("E12").Value = "  -0.43%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.208.52"
$ws.Range("E13").Value = "  -0.74%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.80"
$ws.Range("E14").Value = "  +1.64%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.439.47"
$ws.Range("E15").Value = "  -0.12%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000150"
$ws.Range("E16").Value = "  +0.21%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.727.62"
$ws.Range("E17").Value = "  -0.77%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.56"
$ws.Range("E18").Value = "  +3.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  -1.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.74"
$ws.Range("E20").Value = "  -0.05%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.55"
$ws.Range("E21").Value = "  -2.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.521"
$ws.Range("E23").Value = "  -2.76%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.41"
$ws.Range("E24").Value = "  -0.98%  "

$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.07%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.36"
$ws.Range("E27").Value = "  +0.06%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0902"
$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.96"
$ws.Range("E29").Value = "  +1.96%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.16"
$ws.Range("E30").Value = "  +3.38%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.34"
$ws.Range("E31").Value = "  +12.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "165.60"
$ws.Range("E32").Value = "  -2.00%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.89"
$ws.Range("E33").Value = "  +1.22%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.00"
$ws.Range("E34").Value = "  -0.57%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  +3.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  +0.44%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.974"
$ws.Range("E38").Value = "  -0.06%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "344.73"
$ws.Range("E39").Value = "  +6.26%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.29"
$ws.Range("E40").Value = "  +2.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.08"
$ws.Range("E41").Value = "  -0.93%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "38.53"
$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.82"
$ws.Range("E43").Value = "  +2.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.00"
$ws.Range("E44").Value = "  -0.96%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0582"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0250"
$ws.Range("E47").Value = "  -0.98%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0998"
$ws.Range("E48").Value = "  -0.67%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.998"
$ws.Range("E49").Value = "  -0.12%  "

$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.57"
$ws.Range("E50").Value = "  -2.35%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.08"
$ws.Range("E51").Value = "  +0.46%  "
